# 1. Insert a new column before column B to hold the new 'status_label' field.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("B:B").Insert()

# 2. Header for the new column.
$ws.Cells.Item(1, 2).Value = 'status_label'

# 3. Populate B2:B19 with the French label matching each row's status emoji in column A.
$statusLabels = @{
    '🟥' = 'rouge'
    '🟧' = 'orange'
}
for ($r = 2; $r -le 19; $r++) {
    $statusEmoji = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $statusLabels[$statusEmoji]
}

# 4. The underlying source data was re-pulled/re-ordered for the 2016 / 2017 / 2020 cohorts;
#    re-write those data rows in full (all columns) to match the refreshed row order.
# Row 6: NCT02888964
$ws.Cells.Item(6, 1).Value = '🟥'
$ws.Cells.Item(6, 2).Value = 'rouge'
$ws.Cells.Item(6, 3).Value = 'NCT02888964'
$ws.Cells.Item(6, 4).Value = ''
$ws.Cells.Item(6, 5).Value = '2016'
$ws.Cells.Item(6, 6).Value = 'A Study to Assess Efficacy and Safety of Pioglitazone as Add-On Therapy to Imatinib Mesylate in CP-CML Patients in Major Molecular Response'
$ws.Cells.Item(6, 7).Value = 'ACTIM'
$ws.Cells.Item(6, 8).Value = $false
$ws.Cells.Item(6, 9).Value = $false
$ws.Cells.Item(6, 10).Value = $false

# Row 7: NCT02896842
$ws.Cells.Item(7, 1).Value = '🟥'
$ws.Cells.Item(7, 2).Value = 'rouge'
$ws.Cells.Item(7, 3).Value = 'NCT02896842'
$ws.Cells.Item(7, 4).Value = ''
$ws.Cells.Item(7, 5).Value = '2016'
$ws.Cells.Item(7, 6).Value = 'A Prospective Randomized Phase II Study Evaluating the Monitoring of Imatinib Mesylate (Glivec®) Plasmatic Through Level in Patients Newly Diagnosed With Chronic Phase Chronic Myelogenous Leukaemia (CP-CML).'
$ws.Cells.Item(7, 7).Value = 'OPTIMIMATINIB'
$ws.Cells.Item(7, 8).Value = $false
$ws.Cells.Item(7, 9).Value = $false
$ws.Cells.Item(7, 10).Value = $false

# Row 8: NCT02888990
$ws.Cells.Item(8, 1).Value = '🟥'
$ws.Cells.Item(8, 2).Value = 'rouge'
$ws.Cells.Item(8, 3).Value = 'NCT02888990'
$ws.Cells.Item(8, 4).Value = ''
$ws.Cells.Item(8, 5).Value = '2016'
$ws.Cells.Item(8, 6).Value = 'An Open Label Phase II Study to Evaluate the Efficacy and Safety of Induction and Consolidation Therapy With Dasatinib in Combination With Chemotherapy in Patients Aged 55 Years and Over With Philadelphia Chromosome Positive (Ph+ or BCR-ABL+) Acute Lymphoblastic Leukemia (ALL).'
$ws.Cells.Item(8, 7).Value = 'EWALLPH01'
$ws.Cells.Item(8, 8).Value = $false
$ws.Cells.Item(8, 9).Value = $false
$ws.Cells.Item(8, 10).Value = $false

# Row 9: NCT02883959
$ws.Cells.Item(9, 1).Value = '🟧'
$ws.Cells.Item(9, 2).Value = 'orange'
$ws.Cells.Item(9, 3).Value = 'NCT02883959'
$ws.Cells.Item(9, 4).Value = ''
$ws.Cells.Item(9, 5).Value = '2016'
$ws.Cells.Item(9, 6).Value = ''
$ws.Cells.Item(9, 7).Value = 'Painkiller'
$ws.Cells.Item(9, 8).Value = $false
$ws.Cells.Item(9, 9).Value = $true
$ws.Cells.Item(9, 10).Value = $true

# Row 10: NCT01946750
$ws.Cells.Item(10, 1).Value = '🟥'
$ws.Cells.Item(10, 2).Value = 'rouge'
$ws.Cells.Item(10, 3).Value = 'NCT01946750'
$ws.Cells.Item(10, 4).Value = ''
$ws.Cells.Item(10, 5).Value = '2017'
$ws.Cells.Item(10, 6).Value = ''
$ws.Cells.Item(10, 7).Value = 'SERODIFF'
$ws.Cells.Item(10, 8).Value = $false
$ws.Cells.Item(10, 9).Value = $false
$ws.Cells.Item(10, 10).Value = $false

# Row 13: NCT03115242
$ws.Cells.Item(13, 1).Value = '🟥'
$ws.Cells.Item(13, 2).Value = 'rouge'
$ws.Cells.Item(13, 3).Value = 'NCT03115242'
$ws.Cells.Item(13, 4).Value = ''
$ws.Cells.Item(13, 5).Value = '2017'
$ws.Cells.Item(13, 6).Value = 'Contrast Enhanced Ultrasound of Carotid Plaque in Acute Ischemic Stroke'
$ws.Cells.Item(13, 7).Value = 'CUSCAS'
$ws.Cells.Item(13, 8).Value = $false
$ws.Cells.Item(13, 9).Value = $false
$ws.Cells.Item(13, 10).Value = $false

# Row 15: NCT02894177
$ws.Cells.Item(15, 1).Value = '🟥'
$ws.Cells.Item(15, 2).Value = 'rouge'
$ws.Cells.Item(15, 3).Value = 'NCT02894177'
$ws.Cells.Item(15, 4).Value = ''
$ws.Cells.Item(15, 5).Value = '2020'
$ws.Cells.Item(15, 6).Value = 'Is Transcutaneous Carbon Dioxide Pressure (tcPCO2) Monitoring During Spontaneous Breathing Trials Useful to Predict Extubation Failure in Mechanically Ventilated Patients in the ICU?'
$ws.Cells.Item(15, 7).Value = 'tcPCO2'
$ws.Cells.Item(15, 8).Value = $false
$ws.Cells.Item(15, 9).Value = $false
$ws.Cells.Item(15, 10).Value = $false

# Row 16: NCT03030482
$ws.Cells.Item(16, 1).Value = '🟥'
$ws.Cells.Item(16, 2).Value = 'rouge'
$ws.Cells.Item(16, 3).Value = 'NCT03030482'
$ws.Cells.Item(16, 4).Value = ''
$ws.Cells.Item(16, 5).Value = '2020'
$ws.Cells.Item(16, 6).Value = 'Evaluation of Touch Massage on Anxiety in Critically Ill Patients : a Randomised Controlled Trial Study (REaLAX)'
$ws.Cells.Item(16, 7).Value = 'REaLAX'
$ws.Cells.Item(16, 8).Value = $false
$ws.Cells.Item(16, 9).Value = $false
$ws.Cells.Item(16, 10).Value = $false

